# Refresh computed Leve profit columns (H:N) across the job sheets with
# updated Universalis market-price snapshots (scheduled data-sync runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 11972.625
$ws.Range("I111").Value = 6769.8
$ws.Range("J111").Value = 20644
$ws.Range("K111").Value = 20309.4
$ws.Range("L111").Value = 61932
$ws.Range("M111").Value = -17242.4
$ws.Range("N111").Value = -68066
$ws.Range("H116").Value = 5284.077
$ws.Range("I116").Value = 2849.5
$ws.Range("J116").Value = 7370.857
$ws.Range("K116").Value = 2849.5
$ws.Range("L116").Value = 7370.857
$ws.Range("M116").Value = 592.5
$ws.Range("N116").Value = -14254.857
$ws.Range("H138").Value = 1304
$ws.Range("I138").Value = 1023.4483
$ws.Range("J138").Value = 2660
$ws.Range("K138").Value = 3070.3449
$ws.Range("L138").Value = 7980
$ws.Range("M138").Value = 2069.6551
$ws.Range("N138").Value = -18260

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2909.7058
$ws.Range("I2").Value = 1810.2
$ws.Range("J2").Value = 4480.4287
$ws.Range("K2").Value = 1810.2
$ws.Range("L2").Value = 4480.4287
$ws.Range("M2").Value = -1697.2
$ws.Range("N2").Value = -4706.4287
$ws.Range("H45").Value = 4340.2856
$ws.Range("I45").Value = 3889.5
$ws.Range("J45").Value = 5151.7
$ws.Range("K45").Value = 3889.5
$ws.Range("L45").Value = 5151.7
$ws.Range("M45").Value = -3512.5
$ws.Range("N45").Value = -5905.7
$ws.Range("H61").Value = 2619.4375
$ws.Range("J61").Value = 4577.8184
$ws.Range("L61").Value = 4577.8184
$ws.Range("N61").Value = -5001.8184
$ws.Range("H74").Value = 1567.1818
$ws.Range("I74").Value = 1489.1538
$ws.Range("J74").Value = 1857
$ws.Range("K74").Value = 1489.1538
$ws.Range("L74").Value = 1857
$ws.Range("M74").Value = -615.1538
$ws.Range("N74").Value = -3605
$ws.Range("H77").Value = 1567.1818
$ws.Range("I77").Value = 1489.1538
$ws.Range("J77").Value = 1857
$ws.Range("K77").Value = 7445.769
$ws.Range("L77").Value = 9285
$ws.Range("M77").Value = -3077.769
$ws.Range("N77").Value = -18021
$ws.Range("H116").Value = 2909.7058
$ws.Range("I116").Value = 1810.2
$ws.Range("J116").Value = 4480.4287
$ws.Range("K116").Value = 1810.2
$ws.Range("L116").Value = 4480.4287
$ws.Range("M116").Value = 483.8
$ws.Range("N116").Value = -9068.4287
$ws.Range("H136").Value = 2619.4375
$ws.Range("J136").Value = 4577.8184
$ws.Range("L136").Value = 13733.4552
$ws.Range("N136").Value = -18833.4552

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2909.7058
$ws.Range("I3").Value = 1810.2
$ws.Range("J3").Value = 4480.4287
$ws.Range("K3").Value = 1810.2
$ws.Range("L3").Value = 4480.4287
$ws.Range("M3").Value = -1696.2
$ws.Range("N3").Value = -4708.4287
$ws.Range("H134").Value = 6326.6665
$ws.Range("I134").Value = 3457.7334
$ws.Range("J134").Value = 8375.904
$ws.Range("K134").Value = 10373.2002
$ws.Range("L134").Value = 25127.712
$ws.Range("M134").Value = -7838.200199999999
$ws.Range("N134").Value = -30197.712

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2080.5557
$ws.Range("I31").Value = 1555.921
$ws.Range("K31").Value = 1555.921
$ws.Range("M31").Value = -1260.921
$ws.Range("H34").Value = 2080.5557
$ws.Range("I34").Value = 1555.921
$ws.Range("K34").Value = 1555.921
$ws.Range("M34").Value = -1353.921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 367.8
$ws.Range("I15").Value = 367.8
$ws.Range("K15").Value = 1103.4
$ws.Range("M15").Value = -963.4000000000001
$ws.Range("H18").Value = 630.2353000000001
$ws.Range("I18").Value = 427.2143
$ws.Range("J18").Value = 1577.6666
$ws.Range("K18").Value = 1281.6429
$ws.Range("L18").Value = 4732.9998
$ws.Range("M18").Value = -1112.6429
$ws.Range("N18").Value = -5070.9998
$ws.Range("H92").Value = 1054.3334
$ws.Range("I92").Value = 760.2
$ws.Range("J92").Value = 1201.4
$ws.Range("K92").Value = 2280.6
$ws.Range("L92").Value = 3604.2
$ws.Range("M92").Value = -1032.6
$ws.Range("N92").Value = -6100.200000000001
$ws.Range("H102").Value = 4000
$ws.Range("I102").Value = 4000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 12000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -9566
$ws.Range("N102").ClearContents()
$ws.Range("H107").Value = 41666910
$ws.Range("I107").Value = 166666830
$ws.Range("J107").Value = 275.55554
$ws.Range("K107").Value = 500000490
$ws.Range("L107").Value = 826.66662
$ws.Range("M107").Value = -499998570
$ws.Range("N107").Value = -4666.66662
$ws.Range("H108").Value = 2501.75
$ws.Range("I108").Value = 780.1111
$ws.Range("J108").Value = 7666.6665
$ws.Range("K108").Value = 2340.3333
$ws.Range("L108").Value = 22999.9995
$ws.Range("M108").Value = 539.6667000000002
$ws.Range("N108").Value = -28759.9995
$ws.Range("H129").Value = 2691.25
$ws.Range("I129").Value = 1585
$ws.Range("J129").Value = 2849.2856
$ws.Range("K129").Value = 4755
$ws.Range("L129").Value = 8547.856800000001
$ws.Range("M129").Value = 245
$ws.Range("N129").Value = -18547.8568
$ws.Range("H131").Value = 937.2963
$ws.Range("I131").Value = 540.75
$ws.Range("J131").Value = 1104.2632
$ws.Range("K131").Value = 1622.25
$ws.Range("L131").Value = 3312.7896
$ws.Range("M131").Value = 3417.75
$ws.Range("N131").Value = -13392.7896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 36048.355
$ws.Range("I55").Value = 348.75
$ws.Range("J55").Value = 62823.062
$ws.Range("K55").Value = 348.75
$ws.Range("L55").Value = 62823.062
$ws.Range("M55").Value = -175.75
$ws.Range("N55").Value = -63169.062
$ws.Range("H68").Value = 2749.4736
$ws.Range("I68").Value = 2368.5715
$ws.Range("K68").Value = 2368.5715
$ws.Range("M68").Value = -1619.5715
$ws.Range("H71").Value = 2749.4736
$ws.Range("I71").Value = 2368.5715
$ws.Range("K71").Value = 11842.8575
$ws.Range("M71").Value = -8098.8575
$ws.Range("H136").Value = 1049.5
$ws.Range("I136").Value = 1319.8
$ws.Range("J136").Value = 599
$ws.Range("K136").Value = 3959.4
$ws.Range("L136").Value = 1797
$ws.Range("M136").Value = -1409.4
$ws.Range("N136").Value = -6897

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 619.17645
$ws.Range("I107").Value = 417.3846
$ws.Range("K107").Value = 1252.1538
$ws.Range("M107").Value = 667.8462
